$d = $word.ActiveDocument

# The document body is Java source code formatted at 8pt (sz/szCs=16
# half-points). The edit shrinks that font to 7pt (sz/szCs=14) across the
# whole code block (paragraphs 1-10), leaving the trailing blank
# paragraphs (which use a different, unrelated 10pt style) untouched.

$codeParaCount = 10

for ($i = 1; $i -le $codeParaCount; $i++) {
    $p = $d.Paragraphs.Item($i)
    $r = $p.Range

    if ($r.End -eq $r.Start + 1) {
        # Paragraph consists solely of its end-of-paragraph mark (no run
        # text). Directly assigning Font.Size on such a zero-content
        # range does not stick, so give it a temporary character to
        # carry the formatting, apply the size to paragraph-mark +
        # temp char, then remove the temp char again.
        $r.InsertBefore("X")
        $p2 = $d.Paragraphs.Item($i)
        $p2.Range.Font.Size = 7
        $p2.Range.Font.SizeBi = 7
        $tmp = $d.Range($p2.Range.Start, $p2.Range.Start + 1)
        $tmp.Delete()
    } else {
        $r.Font.Size = 7
        $r.Font.SizeBi = 7
    }
}
